$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = 0.0233936902987173
$ws.Range("B2").Value = 0.02435509861631267
$ws.Range("C2").Value = 0.01638139577444667
$ws.Range("D2").Value = 0.0215334684803855
$ws.Range("E2").Value = 0.02271380992234233
$ws.Range("F2").Value = 0.02396259519812947
$ws.Range("G2").Value = 0.02267511898940877
$ws.Range("H2").Value = 0.0168494776383942
$ws.Range("I2").Value = 0.0057932781162695
$ws.Range("J2").Value = 0.005679314112165876
$ws.Range("K2").Value = 0.02230170902716915
